$d = $word.ActiveDocument

# 1. Fix the typo "Giovanny" -> "Giovanni"
$d.Content.Find.Execute("Giovanny ", $true, $false, $false, $false, $false, $true, 1, $false, "Giovanni ", 2)

# 2. Locate the paragraph that now starts with "Giovanni " (the Boccacio entry) and
#    turn it into a heading (Nadpis4), which drops the bullet/list formatting.
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "Giovanni *") {
        $p.Style = "Nadpis4"
        break
    }
}

# 3. Underline the bold "Boccacio" run.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Boccacio"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
while ($rng.Find.Execute()) {
    $rng.Font.Underline = 1
    $rng.Collapse(0)
}
